$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   C = 249332; E = 1036484754 },
    @{ Row = 53;  C = 141684; E = 590072732 },
    @{ Row = 57;  C = 3714;   E = 138524268 },
    @{ Row = 92;  C = 409227; E = 1596677774 },
    @{ Row = 93;  C = 209631; E = 1309735665 },
    @{ Row = 95;  C = 50797;  E = 933855262 },
    @{ Row = 96;  C = 17309;  E = 795853177 },
    @{ Row = 104; C = 135254; E = 272256995 },
    @{ Row = 110; C = 398;    E = 16733167 },
    @{ Row = 174; C = 226103; E = 900673734 },
    @{ Row = 175; C = 80785;  E = 486184257 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$wb.Save()
